$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '36.743.72'
$ws.Cells.Item(2, 5).Value = '  +3.89%  '

$ws.Cells.Item(3, 4).Value = '1.907.52'
$ws.Cells.Item(3, 5).Value = '  +1.34%  '

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Cells.Item(5, 4).Value = "'248.67"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.84%  '

$ws.Cells.Item(6, 5).Value = '  -0.20%  '

$ws.Cells.Item(7, 5).Value = '  -0.01%  '

$ws.Cells.Item(8, 4).Value = "'46.72"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +7.95%  '

$ws.Cells.Item(9, 4).Value = "'0.372"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +4.55%  '

$ws.Cells.Item(10, 4).Value = "'57.74"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +7.56%  '

$ws.Cells.Item(11, 5).Value = '  +1.28%  '

$ws.Cells.Item(12, 5).Value = '  +2.26%  '

$ws.Cells.Item(13, 4).Value = "'14.65"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +8.45%  '

$ws.Cells.Item(14, 4).Value = "'0.810"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +5.08%  '

$ws.Cells.Item(15, 4).Value = '2.186.77'
$ws.Cells.Item(15, 5).Value = '  +1.35%  '

$ws.Cells.Item(16, 4).Value = "'5.07"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +2.50%  '

$ws.Cells.Item(17, 4).Value = '1.907.30'
$ws.Cells.Item(17, 5).Value = '  +1.07%  '

$ws.Cells.Item(18, 4).Value = '36.752.71'
$ws.Cells.Item(18, 5).Value = '  +3.99%  '

$ws.Cells.Item(19, 4).Value = "'74.16"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +1.00%  '

$ws.Cells.Item(20, 4).Value = '0.0₃0851'
$ws.Cells.Item(20, 5).Value = '  +2.81%  '

$ws.Cells.Item(21, 4).Value = "'13.58"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +5.92%  '

$ws.Cells.Item(22, 4).Value = "'250.10"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +2.20%  '

$ws.Cells.Item(23, 5).Value = '  -1.33%  '

$ws.Cells.Item(24, 5).Value = '  +0.15%  '

$ws.Cells.Item(25, 5).Value = '  -4.06%  '

$ws.Cells.Item(26, 4).Value = "'2.20"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +1.06%  '

$ws.Cells.Item(27, 4).Value = "'166.68"
$ws.Cells.Item(27, 4).Style = 'Normal'

$ws.Cells.Item(28, 4).Value = "'8.73"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.07%  '

$ws.Cells.Item(29, 5).Value = '  +1.79%  '

$ws.Cells.Item(30, 5).Value = '  -0.03%  '

$ws.Cells.Item(31, 5).Value = '  +7.97%  '

$ws.Cells.Item(32, 4).Value = "'0.0612"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +2.66%  '

$ws.Cells.Item(33, 2).Value = 'WEMIXToken'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(33, 4).Value = "'1.93"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +2.92%  '

$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34, 4).Value = "'4.30"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +2.81%  '

$ws.Cells.Item(36, 4).Value = "'0.0868"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +18.21%  '

$ws.Cells.Item(37, 4).Value = "'18.86"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +55.96%  '

$ws.Cells.Item(38, 5).Value = '  -0.98%  '

$ws.Cells.Item(39, 4).Value = "'0.864"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +1.25%  '

$ws.Cells.Item(40, 4).Value = "'1.98"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +1.25%  '

$ws.Cells.Item(41, 4).Value = "'104.17"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +7.10%  '

$ws.Cells.Item(42, 4).Value = "'0.0227"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +4.10%  '

$ws.Cells.Item(43, 4).Value = "'17.70"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +2.30%  '

$ws.Cells.Item(44, 4).Value = "'2.88"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +20.49%  '

$ws.Cells.Item(45, 5).Value = '  +1.10%  '

$ws.Cells.Item(46, 4).Value = '1.346.69'
$ws.Cells.Item(46, 5).Value = '  +2.92%  '

$ws.Cells.Item(47, 4).Value = "'2.36"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -1.40%  '

$ws.Cells.Item(48, 4).Value = "'0.0817"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +0.34%  '

$ws.Cells.Item(49, 4).Value = "'2.82"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +3.15%  '

$ws.Cells.Item(50, 4).Value = "'6.44"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +1.75%  '

$ws.Cells.Item(51, 4).Value = '2.084.71'
$ws.Cells.Item(51, 5).Value = '  +1.12%  '
